$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pre-condition image filenames: forward-slash paths were
# replaced with backslash (Windows-style) paths, and the second file name
# changed from CS+4 to CS-3.
$ws.Range("A2:A6").Value = "PreCondition\CS+3.BMP"
$ws.Range("A7:A11").Value = "PreCondition\CS-3.BMP"

# Extend the selected range from a single cell to A7:A11.
$ws.Range("A7:A11").Select()
